$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.ClearFormats()
}

$ws.Range("D2").Value = "28.100.99"
$ws.Range("E2").Value = "  -0.35%  "
$ws.Range("D3").Value = "1.874.78"
$ws.Range("E3").Value = "  -1.91%  "
$ws.Range("E4").Value = "  +0.15%  "
Set-TextValue "D5" "313.42"
$ws.Range("E5").Value = "  -0.36%  "
$ws.Range("E6").Value = "  +0.13%  "
Set-TextValue "D7" "0.5055"
$ws.Range("E7").Value = "  -0.38%  "
Set-TextValue "D8" "0.3842"
$ws.Range("E8").Value = "  -2.14%  "
Set-TextValue "D9" "0.08583"
$ws.Range("E9").Value = "  -7.85%  "
$ws.Range("E10").Value = "  -2.35%  "
$ws.Range("B11").Value = "OKB"
$ws.Range("C11").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue "D11" "41.36"
$ws.Range("E11").Value = "  -1.46%  "
$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue "D12" "6.315"
$ws.Range("E12").Value = "  -1.34%  "
$ws.Range("B13").Value = "Solana"
$ws.Range("C13").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
Set-TextValue "D13" "20.65"
$ws.Range("E13").Value = "  -1.45%  "
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.873.96"
$ws.Range("E14").Value = "  -2.38%  "
$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextValue "D15" "7.202"
$ws.Range("E15").Value = "  -1.67%  "
$ws.Range("B16").Value = "BinanceUSD"
$ws.Range("C16").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
Set-TextValue "D16" "1.002"
$ws.Range("E16").Value = "  +0.12%  "
$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextValue "D17" "0.00001098"
$ws.Range("E17").Value = "  -2.49%  "
$ws.Range("B18").Value = "Litecoin"
$ws.Range("C18").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextValue "D18" "90.97"
$ws.Range("E18").Value = "  -1.72%  "
$ws.Range("B19").Value = "TRON"
$ws.Range("C19").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
Set-TextValue "D19" "0.06631"
$ws.Range("E19").Value = "  +0.13%  "
$ws.Range("B20").Value = "Avalanche"
$ws.Range("C20").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
Set-TextValue "D20" "18.08"
$ws.Range("E20").Value = "  +0.29%  "
$ws.Range("B21").Value = "Dai"
$ws.Range("C21").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextValue "D21" "1.002"
$ws.Range("E21").Value = "  +0.15%  "
$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
Set-TextValue "D22" "6.101"
$ws.Range("E22").Value = "  -2.09%  "
$ws.Range("B23").Value = "WrappedBTC"
$ws.Range("C23").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D23").Value = "28.136.51"
$ws.Range("E23").Value = "  -0.43%  "
$ws.Range("B24").Value = "Cosmos"
$ws.Range("C24").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue "D24" "11.40"
$ws.Range("E24").Value = "  -1.32%  "
$ws.Range("B25").Value = "Toncoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue "D25" "2.267"
$ws.Range("E25").Value = "  -2.69%  "
$ws.Range("B26").Value = "LidoDAOToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextValue "D26" "2.583"
$ws.Range("E26").Value = "  -0.34%  "
$ws.Range("B27").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C27").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D27").Value = "2.088.54"
$ws.Range("E27").Value = "  -2.46%  "
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue "D28" "20.73"
$ws.Range("E28").Value = "  -2.19%  "
$ws.Range("B29").Value = "Monero"
$ws.Range("C29").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue "D29" "157.09"
$ws.Range("E29").Value = "  -0.59%  "
$ws.Range("B30").Value = "BitcoinCash"
$ws.Range("C30").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-TextValue "D30" "126.33"
$ws.Range("E30").Value = "  -0.75%  "
$ws.Range("B31").Value = "Stellar"
$ws.Range("C31").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue "D31" "0.1057"
$ws.Range("E31").Value = "  -1.74%  "
$ws.Range("B32").Value = "ImmutableX"
$ws.Range("C32").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue "D32" "1.062"
$ws.Range("E32").Value = "  -3.92%  "
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue "D33" "5.616"
$ws.Range("E33").Value = "  -1.06%  "
$ws.Range("B34").Value = "HuobiToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextValue "D34" "3.591"
$ws.Range("E34").Value = "  -0.59%  "
$ws.Range("B35").Value = "FraxShare"
$ws.Range("C35").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue "D35" "9.622"
$ws.Range("E35").Value = "  -0.88%  "
$ws.Range("B36").Value = "VeChain"
$ws.Range("C36").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue "D36" "0.02448"
$ws.Range("E36").Value = "  +0.41%  "
$ws.Range("B37").Value = "Hedera"
$ws.Range("C37").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue "D37" "0.06591"
$ws.Range("E37").Value = "  -1.65%  "
$ws.Range("B38").Value = "Algorand"
$ws.Range("C38").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue "D38" "0.2181"
$ws.Range("E38").Value = "  -1.36%  "
$ws.Range("B39").Value = "ARBITRUM"
$ws.Range("C39").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue "D39" "1.211"
$ws.Range("E39").Value = "  -2.58%  "
$ws.Range("B40").Value = "TrustWalletToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue "D40" "1.244"
$ws.Range("E40").Value = "  -2.92%  "
$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
Set-TextValue "D41" "0.6387"
$ws.Range("E41").Value = "  -2.12%  "
$ws.Range("B42").Value = "Aptos"
$ws.Range("C42").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue "D42" "11.50"
$ws.Range("E42").Value = "  -0.56%  "
$ws.Range("B43").Value = "InternetComputer(DFINITY)"
$ws.Range("C43").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue "D43" "4.898"
$ws.Range("E43").Value = "  -2.64%  "
$ws.Range("B44").Value = "Frax"
$ws.Range("C44").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
Set-TextValue "D44" "1.001"
$ws.Range("E44").Value = "  +0.10%  "
Set-TextValue "D45" "13.27"
$ws.Range("E45").Value = "  -0.61%  "
$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
Set-TextValue "D46" "0.6006"
$ws.Range("E46").Value = "  -2.04%  "
$ws.Range("B47").Value = "WEMIXTOKEN"
$ws.Range("C47").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue "D47" "1.281"
$ws.Range("E47").Value = "  -0.82%  "
$ws.Range("B48").Value = "PancakeSwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue "D48" "3.675"
$ws.Range("E48").Value = "  -1.34%  "
Set-TextValue "D49" "1.991"
$ws.Range("E49").Value = "  -1.88%  "
$ws.Range("B50").Value = "EOS"
$ws.Range("C50").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
Set-TextValue "D50" "1.224"
$ws.Range("E50").Value = "  +2.94%  "
$ws.Range("B51").Value = "Quant"
$ws.Range("C51").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-TextValue "D51" "121.58"
$ws.Range("E51").Value = "  -0.82%  "
